$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" ---
# This value shows up in the Overview sheet (columns "zh-cn"/"de-de", row 2)
# and on the "Status" column of each per-locale sheet, row 2.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width adjustments (report regenerated with narrower date columns) ---
# Overview: columns E ("zh-cn") and F ("de-de")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de sheets: column C ("Status")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
